$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "JOSE VICENTE GUZMAN" worker block (two rows: period 1610 and 1609).
# These are currently rows 18 and 19 in the data table.
$ws.Rows.Item(18).Delete() | Out-Null
$ws.Rows.Item(18).Delete() | Out-Null

# The remaining three workers (IVAN MENDOZA RAMIREZ, LUIS CARLOS GUTIERREZ VEGA,
# JHONATHAN RECUERO MORELO) are re-listed grouped by period: all three for period
# 1609 first (rows 16-18), then all three again for period 1610 (rows 19-21).
$ws.Range("B16:B21").Value = "CC"

$ws.Range("C16").Value = "1129495372"
$ws.Range("D16").Value = "IVAN MENDOZA RAMIREZ"
$ws.Range("E16").Value = "1609"

$ws.Range("C17").Value = "73559861"
$ws.Range("D17").Value = "LUIS CARLOS GUTIERREZ VEGA"
$ws.Range("E17").Value = "1609"

$ws.Range("C18").Value = "73212463"
$ws.Range("D18").Value = "JHONATHAN RECUERO MORELO"
$ws.Range("E18").Value = "1609"

$ws.Range("C19").Value = "1129495372"
$ws.Range("D19").Value = "IVAN MENDOZA RAMIREZ"
$ws.Range("E19").Value = "1610"

$ws.Range("C20").Value = "73559861"
$ws.Range("D20").Value = "LUIS CARLOS GUTIERREZ VEGA"
$ws.Range("E20").Value = "1610"

$ws.Range("C21").Value = "73212463"
$ws.Range("D21").Value = "JHONATHAN RECUERO MORELO"
$ws.Range("E21").Value = "1610"

# Updated totals: worker count 4 -> 3, and the overdue amount total.
$ws.Range("C13").Value = 3
$ws.Range("E11").Value = 165468
